$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 872.25
$ws.Range("I12").Value = 498
$ws.Range("K12").Value = 498
$ws.Range("M12").Value = -328
$ws.Range("H17").Value = 244848.31
$ws.Range("J17").Value = 269193.06
$ws.Range("L17").Value = 807579.1799999999
$ws.Range("N17").Value = -807915.1799999999
$ws.Range("H40").Value = 3650
$ws.Range("H41").Value = 948.2778
$ws.Range("J41").Value = 641.6667
$ws.Range("L41").Value = 641.6667
$ws.Range("N41").Value = -1521.6667
$ws.Range("H61").Value = 498
$ws.Range("I61").Value = 498
$ws.Range("K61").Value = 1494
$ws.Range("M61").Value = -1322
$ws.Range("H101").Value = 1454.9
$ws.Range("I101").Value = 1554.3334
$ws.Range("K101").Value = 4663.0002
$ws.Range("M101").Value = -3041.0002
$ws.Range("H112").Value = 1785.762
$ws.Range("J112").Value = 1868.3889
$ws.Range("L112").Value = 5605.1667
$ws.Range("N112").Value = -7821.1667
$ws.Range("H137").Value = 2056.5557
$ws.Range("I137").Value = 1804
$ws.Range("K137").Value = 5412
$ws.Range("M137").Value = -2862
$ws.Range("H138").Value = 1638.2858
$ws.Range("I138").Value = 1160.7307
$ws.Range("J138").Value = 1920.4773
$ws.Range("K138").Value = 3482.1921
$ws.Range("L138").Value = 5761.4319
$ws.Range("M138").Value = 1657.8079
$ws.Range("N138").Value = -16041.4319

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 21701.4
$ws.Range("I28").Value = 19995
$ws.Range("J28").Value = 22839
$ws.Range("K28").Value = 19995
$ws.Range("L28").Value = 22839
$ws.Range("M28").Value = -19803
$ws.Range("N28").Value = -23223
$ws.Range("H31").Value = 7510.8184
$ws.Range("I31").Value = 5210
$ws.Range("K31").Value = 5210
$ws.Range("M31").Value = -4916
$ws.Range("H32").Value = 4936.381
$ws.Range("I32").Value = 4290.5127
$ws.Range("K32").Value = 4290.5127
$ws.Range("M32").Value = -4003.5127
$ws.Range("H74").Value = 1278.6875
$ws.Range("I74").Value = 1278.6875
$ws.Range("K74").Value = 1278.6875
$ws.Range("M74").Value = -404.6875
$ws.Range("H77").Value = 1278.6875
$ws.Range("I77").Value = 1278.6875
$ws.Range("K77").Value = 6393.4375
$ws.Range("M77").Value = -2025.4375
$ws.Range("H99").Value = 21701.4
$ws.Range("I99").Value = 19995
$ws.Range("J99").Value = 22839
$ws.Range("K99").Value = 19995
$ws.Range("L99").Value = 22839
$ws.Range("M99").Value = -17000
$ws.Range("N99").Value = -28829
$ws.Range("H132").Value = 1619.1029
$ws.Range("I132").Value = 1543.7903
$ws.Range("J132").Value = 2397.3333
$ws.Range("K132").Value = 4631.3709
$ws.Range("L132").Value = 7191.999899999999
$ws.Range("M132").Value = -2101.3709
$ws.Range("N132").Value = -12251.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 699.25
$ws.Range("I7").Value = 510
$ws.Range("J7").Value = 762.3333
$ws.Range("K7").Value = 510
$ws.Range("L7").Value = 762.3333
$ws.Range("M7").Value = -397
$ws.Range("N7").Value = -988.3333
$ws.Range("H94").Value = 1905.3334
$ws.Range("I94").Value = 2045.9
$ws.Range("J94").Value = 1202.5
$ws.Range("K94").Value = 2045.9
$ws.Range("L94").Value = 1202.5
$ws.Range("M94").Value = -1594.9
$ws.Range("N94").Value = -2104.5
$ws.Range("H102").Value = 21407.666
$ws.Range("J102").Value = 30611.5
$ws.Range("L102").Value = 30611.5
$ws.Range("N102").Value = -37101.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10533.629
$ws.Range("I31").Value = 1677.4828
$ws.Range("J31").Value = 53338.332
$ws.Range("K31").Value = 1677.4828
$ws.Range("L31").Value = 53338.332
$ws.Range("M31").Value = -1382.4828
$ws.Range("N31").Value = -53928.332
$ws.Range("H34").Value = 10533.629
$ws.Range("I34").Value = 1677.4828
$ws.Range("J34").Value = 53338.332
$ws.Range("K34").Value = 1677.4828
$ws.Range("L34").Value = 53338.332
$ws.Range("M34").Value = -1475.4828
$ws.Range("N34").Value = -53742.332
$ws.Range("H124").Value = 99999.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 99999.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 99999.5
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -104909.5
$ws.Range("H132").Value = 3960.0952
$ws.Range("I132").Value = 4925.533
$ws.Range("J132").Value = 1546.5
$ws.Range("K132").Value = 14776.599
$ws.Range("L132").Value = 4639.5
$ws.Range("M132").Value = -12246.599
$ws.Range("N132").Value = -9699.5
$ws.Range("H141").Value = 420556
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 420556
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 420556
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -430916

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 349.66666
$ws.Range("I2").Value = 553.9167
$ws.Range("J2").Value = 145.41667
$ws.Range("K2").Value = 3323.5002
$ws.Range("L2").Value = 872.5000200000001
$ws.Range("M2").Value = -3210.5002
$ws.Range("N2").Value = -1098.50002
$ws.Range("H23").Value = 104.21429
$ws.Range("I23").Value = 126.28571
$ws.Range("K23").Value = 378.85713
$ws.Range("M23").Value = -143.85713
$ws.Range("H35").Value = 400
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 1200
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -1776
$ws.Range("H37").Value = 250063730
$ws.Range("J37").Value = 250063730
$ws.Range("L37").Value = 750191190
$ws.Range("N37").Value = -750191414
$ws.Range("H38").Value = 170.42857
$ws.Range("I38").Value = 32
$ws.Range("J38").Value = 1001
$ws.Range("K38").Value = 96
$ws.Range("L38").Value = 3003
$ws.Range("M38").Value = 251
$ws.Range("N38").Value = -3697
$ws.Range("H124").Value = 1500
$ws.Range("I124").Value = 1500
$ws.Range("K124").Value = 4500
$ws.Range("M124").Value = 410
$ws.Range("H126").Value = 765
$ws.Range("I126").Value = 765
$ws.Range("K126").Value = 2295
$ws.Range("M126").Value = 2645
$ws.Range("H131").Value = 6904.579
$ws.Range("I131").Value = 15840
$ws.Range("J131").Value = 1692.25
$ws.Range("K131").Value = 47520
$ws.Range("L131").Value = 5076.75
$ws.Range("M131").Value = -42480
$ws.Range("N131").Value = -15156.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1743.86
$ws.Range("I132").Value = 1752.3636
$ws.Range("J132").Value = 1681.5
$ws.Range("K132").Value = 5257.0908
$ws.Range("L132").Value = 5044.5
$ws.Range("M132").Value = -2727.0908
$ws.Range("N132").Value = -10104.5
$ws.Range("H136").Value = 30553.191
$ws.Range("J136").Value = 30553.191
$ws.Range("L136").Value = 91659.573
$ws.Range("N136").Value = -96759.573

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11479.714
$ws.Range("I7").Value = 14922.526
$ws.Range("J7").Value = 4211.5557
$ws.Range("K7").Value = 14922.526
$ws.Range("L7").Value = 4211.5557
$ws.Range("M7").Value = -14810.526
$ws.Range("N7").Value = -4435.5557
$ws.Range("H68").Value = 2057
$ws.Range("I68").Value = 1949.75
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 1949.75
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -1200.75
$ws.Range("N68").Value = -3698
$ws.Range("H71").Value = 2057
$ws.Range("I71").Value = 1949.75
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 9748.75
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -6004.75
$ws.Range("N71").Value = -18488
$ws.Range("H82").Value = 2462.111
$ws.Range("J82").Value = 2673.25
$ws.Range("L82").Value = 2673.25
$ws.Range("N82").Value = -3395.25
$ws.Range("H85").Value = 2462.111
$ws.Range("J85").Value = 2673.25
$ws.Range("L85").Value = 2673.25
$ws.Range("N85").Value = -5169.25
$ws.Range("H87").Value = 25189
$ws.Range("J87").Value = 25189
$ws.Range("L87").Value = 25189
$ws.Range("N87").Value = -27435
$ws.Range("H90").Value = 25189
$ws.Range("J90").Value = 25189
$ws.Range("L90").Value = 75567
$ws.Range("N90").Value = -86799
$ws.Range("H93").Value = 23510.25
$ws.Range("I93").Value = 3002.7273
$ws.Range("J93").Value = 68626.8
$ws.Range("K93").Value = 3002.7273
$ws.Range("L93").Value = 68626.8
$ws.Range("M93").Value = -1754.7273
$ws.Range("N93").Value = -71122.8
$ws.Range("H126").Value = 11479.714
$ws.Range("I126").Value = 14922.526
$ws.Range("J126").Value = 4211.5557
$ws.Range("K126").Value = 44767.578
$ws.Range("L126").Value = 12634.6671
$ws.Range("M126").Value = -42297.578
$ws.Range("N126").Value = -17574.6671
$ws.Range("H136").Value = 3216.3572
$ws.Range("I136").Value = 2457.6365
$ws.Range("K136").Value = 7372.9095
$ws.Range("M136").Value = -4822.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 7999.75
$ws.Range("I22").Value = 6000
$ws.Range("J22").Value = 9999.5
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 9999.5
$ws.Range("M22").Value = -5707
$ws.Range("N22").Value = -10585.5
$ws.Range("H132").Value = 2653.3704
$ws.Range("I132").Value = 2862.2
$ws.Range("J132").Value = 2056.7144
$ws.Range("K132").Value = 8586.599999999999
$ws.Range("L132").Value = 6170.1432
$ws.Range("M132").Value = -6056.599999999999
$ws.Range("N132").Value = -11230.1432

Write-Host "Applied all Leviathan Profits updates"